# Scheduled-runner update: refresh cached market-board price figures
# (currentAveragePrice / NQ / HQ, leve price, and profit columns) across
# the per-class "Chocobo Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1422.2222
$ws.Range("J40").Value = 1700
$ws.Range("L40").Value = 1700
$ws.Range("N40").Value = -2050

$ws.Range("H51").Value = 5669
$ws.Range("I51").Value = 3566.6667
$ws.Range("J51").Value = 8822.5
$ws.Range("K51").Value = 3566.6667
$ws.Range("L51").Value = 8822.5
$ws.Range("M51").Value = -3082.6667
$ws.Range("N51").Value = -9790.5

$ws.Range("H70").Value = 2400.9092
$ws.Range("I70").Value = 1450
$ws.Range("J70").Value = 3542
$ws.Range("K70").Value = 4350
$ws.Range("L70").Value = 10626
$ws.Range("M70").Value = -4080
$ws.Range("N70").Value = -11166

$ws.Range("H73").Value = 2400.9092
$ws.Range("I73").Value = 1450
$ws.Range("J73").Value = 3542
$ws.Range("K73").Value = 4350
$ws.Range("L73").Value = 10626
$ws.Range("M73").Value = -3414
$ws.Range("N73").Value = -12498

$ws.Range("H75").Value = 29157
$ws.Range("J75").Value = 29157
$ws.Range("L75").Value = 29157
$ws.Range("N75").Value = -31029

$ws.Range("H78").Value = 29157
$ws.Range("J78").Value = 29157
$ws.Range("L78").Value = 87471
$ws.Range("N78").Value = -96831

$ws.Range("H116").Value = 7116.5415
$ws.Range("I116").Value = 2784.125
$ws.Range("J116").Value = 9282.75
$ws.Range("K116").Value = 2784.125
$ws.Range("L116").Value = 9282.75
$ws.Range("M116").Value = 657.875
$ws.Range("N116").Value = -16166.75

$ws.Range("H120").Value = 40761
$ws.Range("J120").Value = 40761
$ws.Range("L120").Value = 40761
$ws.Range("N120").Value = -50437

$ws.Range("H129").Value = 840.37
$ws.Range("I129").Value = 338.2
$ws.Range("J129").Value = 866.8
$ws.Range("K129").Value = 1014.6
$ws.Range("L129").Value = 2600.4
$ws.Range("M129").Value = 3985.4
$ws.Range("N129").Value = -12600.4

$ws.Range("H132").Value = 47626692
$ws.Range("I132").Value = 66675556
$ws.Range("J132").Value = 4532.6665
$ws.Range("K132").Value = 200026668
$ws.Range("L132").Value = 13597.9995
$ws.Range("M132").Value = -200024138
$ws.Range("N132").Value = -18657.9995

$ws.Range("H137").Value = 1222925.2
$ws.Range("I137").Value = 1588669.9
$ws.Range("J137").Value = 3776.5557
$ws.Range("K137").Value = 4766009.699999999
$ws.Range("L137").Value = 11329.6671
$ws.Range("M137").Value = -4763459.699999999
$ws.Range("N137").Value = -16429.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 31804.348
$ws.Range("J109").Value = 31804.348
$ws.Range("L109").Value = 31804.348
$ws.Range("N109").Value = -34578.348

$ws.Range("H132").Value = 2440.5386
$ws.Range("I132").Value = 1432.9231
$ws.Range("J132").Value = 3448.1538
$ws.Range("K132").Value = 4298.7693
$ws.Range("L132").Value = 10344.4614
$ws.Range("M132").Value = -1768.7693
$ws.Range("N132").Value = -15404.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8579.434999999999
$ws.Range("I31").Value = 3135.5454
$ws.Range("K31").Value = 3135.5454
$ws.Range("M31").Value = -2840.5454

$ws.Range("H34").Value = 8579.434999999999
$ws.Range("I34").Value = 3135.5454
$ws.Range("K34").Value = 3135.5454
$ws.Range("M34").Value = -2933.5454

$ws.Range("H62").Value = 3515
$ws.Range("I62").Value = 3002.5
$ws.Range("J62").Value = 3720
$ws.Range("K62").Value = 3002.5
$ws.Range("L62").Value = 3720
$ws.Range("M62").Value = -2378.5
$ws.Range("N62").Value = -4968

$ws.Range("H65").Value = 3515
$ws.Range("I65").Value = 3002.5
$ws.Range("J65").Value = 3720
$ws.Range("K65").Value = 15012.5
$ws.Range("L65").Value = 18600
$ws.Range("M65").Value = -11892.5
$ws.Range("N65").Value = -24840

$ws.Range("H68").Value = 99999
$ws.Range("J68").Value = 99999
$ws.Range("L68").Value = 99999
$ws.Range("N68").Value = -101497

$ws.Range("H71").Value = 99999
$ws.Range("J71").Value = 99999
$ws.Range("L71").Value = 299997
$ws.Range("N71").Value = -307485

$ws.Range("H99").Value = 10003884
$ws.Range("I99").Value = 18183362
$ws.Range("J99").Value = 6744.4443
$ws.Range("K99").Value = 18183362
$ws.Range("L99").Value = 6744.4443
$ws.Range("M99").Value = -18181864
$ws.Range("N99").Value = -9740.444299999999

$ws.Range("H122").Value = 4041.5
$ws.Range("I122").Value = 1916.3334
$ws.Range("J122").Value = 6166.6665
$ws.Range("K122").Value = 5749.0002
$ws.Range("L122").Value = 18499.9995
$ws.Range("M122").Value = -3299.0002
$ws.Range("N122").Value = -23399.9995

$ws.Range("H126").Value = 10003884
$ws.Range("I126").Value = 18183362
$ws.Range("J126").Value = 6744.4443
$ws.Range("K126").Value = 54550086
$ws.Range("L126").Value = 20233.3329
$ws.Range("M126").Value = -54547616
$ws.Range("N126").Value = -25173.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 250001520
$ws.Range("I114").Value = 500000060
$ws.Range("J114").Value = 3000
$ws.Range("K114").Value = 1500000180
$ws.Range("L114").Value = 9000
$ws.Range("M114").Value = -1499996926
$ws.Range("N114").Value = -15508

$ws.Range("H117").Value = 1598.7778
$ws.Range("J117").Value = 2026.6666
$ws.Range("L117").Value = 6079.9998
$ws.Range("N117").Value = -12963.9998

$ws.Range("H129").Value = 3296.9167
$ws.Range("I129").Value = 3621.6667
$ws.Range("J129").Value = 2972.1667
$ws.Range("K129").Value = 10865.0001
$ws.Range("L129").Value = 8916.500100000001
$ws.Range("M129").Value = -5865.000100000001
$ws.Range("N129").Value = -18916.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 30735.889
$ws.Range("J57").Value = 30735.889
$ws.Range("L57").Value = 30735.889
$ws.Range("N57").Value = -32375.889

$ws.Range("H80").Value = 27782100
$ws.Range("I80").Value = 41671250
$ws.Range("J80").Value = 3800
$ws.Range("K80").Value = 41671250
$ws.Range("L80").Value = 3800
$ws.Range("M80").Value = -41670252
$ws.Range("N80").Value = -5796

$ws.Range("H83").Value = 27782100
$ws.Range("I83").Value = 41671250
$ws.Range("J83").Value = 3800
$ws.Range("K83").Value = 208356250
$ws.Range("L83").Value = 19000
$ws.Range("M83").Value = -208351258
$ws.Range("N83").Value = -28984

$ws.Range("H96").Value = 31790.25
$ws.Range("J96").Value = 31790.25
$ws.Range("L96").Value = 31790.25
$ws.Range("N96").Value = -37282.25

$ws.Range("H113").Value = 2887.5
$ws.Range("I113").Value = 2883.3333
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 2883.3333
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -713.3332999999998
$ws.Range("N113").Value = -7240

$ws.Range("H122").Value = 3933.625
$ws.Range("I122").Value = 2955.818
$ws.Range("J122").Value = 6084.8
$ws.Range("K122").Value = 8867.454000000002
$ws.Range("L122").Value = 18254.4
$ws.Range("M122").Value = -6417.454000000002
$ws.Range("N122").Value = -23154.4

$ws.Range("H141").Value = 39500
$ws.Range("I141").Value = 36000
$ws.Range("J141").Value = 43000
$ws.Range("K141").Value = 36000
$ws.Range("L141").Value = 43000
$ws.Range("M141").Value = -30820
$ws.Range("N141").Value = -53360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 15000
$ws.Range("J40").Value = 8333.333000000001
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 8333.333000000001
$ws.Range("M40").Value = -14864
$ws.Range("N40").Value = -8605.333000000001

$ws.Range("H122").Value = 7825
$ws.Range("I122").Value = 4800
$ws.Range("K122").Value = 14400
$ws.Range("M122").Value = -11950

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H140").Value = 72009.55499999999
$ws.Range("J140").Value = 72009.55499999999
$ws.Range("L140").Value = 72009.55499999999
$ws.Range("N140").Value = -82369.55499999999

$ws.Range("H141").Value = 32247.5
$ws.Range("J141").Value = 32247.5
$ws.Range("L141").Value = 32247.5
$ws.Range("N141").Value = -42607.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 52562.363
$ws.Range("J46").Value = 52562.363
$ws.Range("L46").Value = 52562.363
$ws.Range("N46").Value = -53024.363

$ws.Range("H57").Value = 19466.666
$ws.Range("J57").Value = 19466.666
$ws.Range("L57").Value = 19466.666
$ws.Range("N57").Value = -20974.666

$ws.Range("H134").Value = 52562.363
$ws.Range("J134").Value = 52562.363
$ws.Range("L134").Value = 157687.089
$ws.Range("N134").Value = -162757.089

$ws.Range("H140").Value = 35427.4
$ws.Range("J140").Value = 35427.4
$ws.Range("L140").Value = 35427.4
$ws.Range("N140").Value = -45787.4

$ws.Range("H141").Value = 37706.25
$ws.Range("J141").Value = 37706.25
$ws.Range("L141").Value = 37706.25
$ws.Range("N141").Value = -48066.25
